$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.357.83'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.35%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.847.63'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.26%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.40'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.21%  '

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.42%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.10%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07583'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.12%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2908'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.15%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.61'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.29%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07749'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.07%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.847.85'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.37%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.019'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.33%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6796'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.31%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001041'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.76%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.09'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.69%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.116'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.91%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.349.92'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.43%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '229.16'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.07%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.34'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.16%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.07%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.429'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.40%  '

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.05%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '159.01'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.26%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1392'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.52%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.431'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.24%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.66'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.39%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.432'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.84%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.468'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.12%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05667'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.42%  '

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.61%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.041'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.21%  '

# Row 33
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.156'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.91%  '

# Row 34
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.822'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.62%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7005'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.14%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.581'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.23%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01831'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.97%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.234.91'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.27%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.13%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.386'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.65%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8972'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.44%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9995'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.20%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.25'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.49%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.46'

# Row 45
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000119'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.49%  '

# Row 46
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.128'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.17%  '

# Row 47
$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.3997'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.61%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1152'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.22%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.969'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.51%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.675'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.61%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05697'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.39%  '
